$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 238; existing rows 238-255 shift down to 240-257.
$ws.Rows("238:239").Insert()

# New row 238 (week of 2022-06-02, quality "Primera")
$ws.Range("A238").Value = 4
$ws.Range("B238").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C238").Value = "Los Lagos"
$ws.Range("D238").Value = 44714
$ws.Range("E238").Value = 10
$ws.Range("F238").Value = "Fruta"
$ws.Range("G238").Value = 100104
$ws.Range("H238").Value = "Frutos de pepita"
$ws.Range("I238").Value = 100104005
$ws.Range("J238").Value = "Pera"
$ws.Range("K238").Value = "Packham's Triumph"
$ws.Range("L238").Value = "Primera"
$ws.Range("M238").Value = 400
$ws.Range("N238").Value = 15000
$ws.Range("O238").Value = 16000
$ws.Range("P238").Value = 15500
$ws.Range("Q238").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R238").Value = "Región de O'Higgins"
$ws.Range("S238").Value = 1033
$ws.Range("T238").Value = 15

# New row 239 (week of 2022-06-02, quality "Segunda")
$ws.Range("A239").Value = 4
$ws.Range("B239").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C239").Value = "Los Lagos"
$ws.Range("D239").Value = 44714
$ws.Range("E239").Value = 10
$ws.Range("F239").Value = "Fruta"
$ws.Range("G239").Value = 100104
$ws.Range("H239").Value = "Frutos de pepita"
$ws.Range("I239").Value = 100104005
$ws.Range("J239").Value = "Pera"
$ws.Range("K239").Value = "Packham's Triumph"
$ws.Range("L239").Value = "Segunda"
$ws.Range("M239").Value = 200
$ws.Range("N239").Value = 13000
$ws.Range("O239").Value = 13000
$ws.Range("P239").Value = 13000
$ws.Range("Q239").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R239").Value = "Región de O'Higgins"
$ws.Range("S239").Value = 867
$ws.Range("T239").Value = 15
